$data = @(
    @(4.229360163211825, -4.647700071334836, -1.350979804992678),
    @(4.624738931655887, -3.006192684173583, -2.288825809955597),
    @(-6.357042789459214, -8.575422286987301, 0.2564473152160626),
    @(0.5462948679924189, -6.825089752674098, -0.5449948012828849),
    @(0.9978208541870109, -3.706368923187252, -1.51154860854149),
    @(0.2441467046737658, -3.069634318351747, -2.925750926136973),
    @(1.612907171249393, -4.853008508682255, -1.383459806442257),
    @(-1.189411103725463, -6.66196793317795, 2.149218022823342),
    @(-6.306459784507723, -6.704558491706845, 4.024554014205929),
    @(1.663261890411397, -5.016231019049863, 3.32933139801026),
    @(6.837078571319559, 0.09169325232505482, 5.010437965393057),
    @(2.122651159763334, -0.6034613586962226, 2.880795598030089),
    @(-0.7601926326751736, 2.327319413423542, 5.726811170578007),
    @(4.228423535823836, -2.42055988311769, 5.180934607982632),
    @(2.933720350265484, -5.437817335128782, 5.24878549575806),
    @(-2.91136687994004, -3.510188579559316, 5.247701197862618),
    @(-4.280053377151489, -0.6633338928222658, 4.682214915752411),
    @(-2.581492483615869, 0.5014263689517995, 6.192452192306466),
    @(-0.986172676086416, 1.025731801986666, -7.288565635681123),
    @(1.218793094158184, -5.294871598482164, -0.7042694091796555),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-excess last data row (original sheet had rows 2-22; final has 2-21)
$ws.Rows.Item(22).Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Host "Updated range:" $ws.UsedRange.Address()
